# Auto-generated edit script applying per-cell value updates from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.905.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.300.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.657.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.312.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.880.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "25.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.81"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.64%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("E35").Value = "  -3.24%  "
$ws.Range("E36").Value = "  -1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0685"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.999.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.522.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.67%  "
